# Applies the "Removed CTI and DWI input, updated the format and set default
# parameters to analyse test_data" edit described by the diff.
#
# Color constants (Excel COM colors are packed as R + G*256 + B*65536):
#   red   (FFFF0000) -> "off"/"no" conditional formatting font color = 255
#   green (FF00B050) -> "on"/"yes" conditional formatting font color = 5287936

$wb = $excel.ActiveWorkbook

$RED = 255
$GREEN = 5287936

# ---------------------------------------------------------------------------
# Sheet 1: activity_analysis
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("activity_analysis")

# Default parameters: turn these four switches from "on" to "off"
$ws1.Range("B2").Value2 = "off"   # erp.sensor_enable
$ws1.Range("B9").Value2 = "off"   # ers_erd.sensor_enable
$ws1.Range("B10").Value2 = "off"  # ers_erd.roi_enable
$ws1.Range("B12").Value2 = "off"  # ers_erd.mapping_enable

# Recolor the existing conditional formatting: off -> red, on -> green
$fcs1 = $ws1.Range("B1:B1048576").FormatConditions
for ($i = 1; $i -le $fcs1.Count; $i++) {
    $fc = $fcs1.Item($i)
    if ($fc.Formula1 -eq '="off"') {
        $fc.Font.Color = $RED
    } elseif ($fc.Formula1 -eq '="on"') {
        $fc.Font.Color = $GREEN
    }
}

# Move the saved selection
$ws1.Activate()
$ws1.Range("B13").Select()

# ---------------------------------------------------------------------------
# Sheet 2: connectivity_analysis
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("connectivity_analysis")

# Default parameters: turn ica_conn.enable from "off" to "on"
$ws2.Range("B2").Value2 = "on"

# Remove the separate B11-only cellIs formatting (erp/triggers row no longer
# gets its own on/off rule; it is absorbed into the column-wide rule below).
$fcsB11 = $ws2.Range("B11").FormatConditions
while ($fcsB11.Count -gt 0) {
    $fcsB11.Item(1).Delete()
}

# Extend the beginsWith on/off rule (previously split across B2:B10 & B12:B20)
# to the whole column, and recolor off -> red, on -> green.
$fullRange2 = $ws2.Range("B1:B1048576")
$fcsMulti = $ws2.Range("B2:B10,B12:B20").FormatConditions
for ($i = 1; $i -le $fcsMulti.Count; $i++) {
    $fc = $fcsMulti.Item($i)
    $isOff = $fc.Formula1 -like '*"off"*'
    $isOn = $fc.Formula1 -like '*"on"*'
    $fc.ModifyAppliesToRange($fullRange2)
    if ($isOff) {
        $fc.Formula1 = '=LEFT(B1,LEN("off"))="off"'
        $fc.Font.Color = $RED
        $fc.Priority = 1
    } elseif ($isOn) {
        $fc.Formula1 = '=LEFT(B1,LEN("on"))="on"'
        $fc.Font.Color = $GREEN
        $fc.Priority = 2
    }
}

# Move the saved selection
$ws2.Activate()
$ws2.Range("B3").Select()

# ---------------------------------------------------------------------------
# Sheet 3: statistical_analysis
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("statistical_analysis")

# Default parameters: turn stats.demean from "no" to "yes"
$ws3.Range("B4").Value2 = "yes"

# Extend the containsText no/yes rule from B2:B12 to the whole column, and
# recolor no -> red, yes -> green.
$fullRange3 = $ws3.Range("B1:B1048576")
$fcs3 = $ws3.Range("B2:B12").FormatConditions
for ($i = 1; $i -le $fcs3.Count; $i++) {
    $fc = $fcs3.Item($i)
    $isNo = $fc.Formula1 -like '*"no"*'
    $isYes = $fc.Formula1 -like '*"yes"*'
    $fc.ModifyAppliesToRange($fullRange3)
    if ($isNo) {
        $fc.Formula1 = '=NOT(ISERROR(SEARCH("no",B1)))'
        $fc.Font.Color = $RED
    } elseif ($isYes) {
        $fc.Formula1 = '=NOT(ISERROR(SEARCH("yes",B1)))'
        $fc.Font.Color = $GREEN
    }
}

# Move the saved selection
$ws3.Activate()
$ws3.Range("B5").Select()
